$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old first column (A) held redundant values that are now dropped: the
# remaining columns (old B:F) shift one position to the left (new A:E).
$ws.Columns("A").Delete()

# Fix the "MODEL_CONDITION" header text (now in column D) -> "MODELCONDITION".
$ws.Range("D1").Value = "MODELCONDITION"
